$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.267161011695862
$ws.Range("B1").Value = 1.682010293006897
$ws.Range("C1").Value = 1.421332597732544
$ws.Range("D1").Value = 2.067212343215942
$ws.Range("E1").Value = 3.324020624160767
